$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 50003976
$ws.Range("I62").Value = 66668732
$ws.Range("J62").Value = 9700
$ws.Range("K62").Value = 66668732
$ws.Range("L62").Value = 9700
$ws.Range("M62").Value = -66668108
$ws.Range("N62").Value = -10948

$ws.Range("H65").Value = 50003976
$ws.Range("I65").Value = 66668732
$ws.Range("J65").Value = 9700
$ws.Range("K65").Value = 333343660
$ws.Range("L65").Value = 48500
$ws.Range("M65").Value = -333340540
$ws.Range("N65").Value = -54740

$ws.Range("H138").Value = 5749493.5
$ws.Range("J138").Value = 8931692
$ws.Range("L138").Value = 26795076
$ws.Range("N138").Value = -26805356

$ws.Range("H141").Value = 5000
$ws.Range("I141").Value = 5000
$ws.Range("K141").Value = 15000
$ws.Range("M141").Value = -9820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9210.691999999999
$ws.Range("J32").Value = 25450.637
$ws.Range("L32").Value = 25450.637
$ws.Range("N32").Value = -26024.637

$ws.Range("H45").Value = 14363.5
$ws.Range("J45").Value = 2050
$ws.Range("L45").Value = 2050
$ws.Range("N45").Value = -2804

$ws.Range("H61").Value = 3411.4783
$ws.Range("I61").Value = 2291.1304
$ws.Range("K61").Value = 2291.1304
$ws.Range("M61").Value = -2079.1304

$ws.Range("H74").Value = 9172.280000000001
$ws.Range("I74").Value = 1916.7368
$ws.Range("K74").Value = 1916.7368
$ws.Range("M74").Value = -1042.7368

$ws.Range("H77").Value = 9172.280000000001
$ws.Range("I77").Value = 1916.7368
$ws.Range("K77").Value = 9583.683999999999
$ws.Range("M77").Value = -5215.683999999999

$ws.Range("H97").Value = 1411.4231
$ws.Range("I97").Value = 1008.9545
$ws.Range("J97").Value = 3625
$ws.Range("K97").Value = 1008.9545
$ws.Range("L97").Value = 3625
$ws.Range("M97").Value = -512.9545000000001
$ws.Range("N97").Value = -4617

$ws.Range("H132").Value = 2651.4167
$ws.Range("I132").Value = 2483.853
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 7451.559
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -4921.559
$ws.Range("N132").Value = -21560

$ws.Range("H136").Value = 3411.4783
$ws.Range("I136").Value = 2291.1304
$ws.Range("K136").Value = 6873.3912
$ws.Range("M136").Value = -4323.3912

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H80").Value = 360.8889
$ws.Range("I80").Value = 444.2
$ws.Range("J80").Value = 328.84616
$ws.Range("K80").Value = 444.2
$ws.Range("L80").Value = 328.84616
$ws.Range("M80").Value = 553.8
$ws.Range("N80").Value = -2324.84616

$ws.Range("H83").Value = 360.8889
$ws.Range("I83").Value = 444.2
$ws.Range("J83").Value = 328.84616
$ws.Range("K83").Value = 2221
$ws.Range("L83").Value = 1644.2308
$ws.Range("M83").Value = 2771
$ws.Range("N83").Value = -11628.2308

$ws.Range("H94").Value = 958.73914
$ws.Range("I94").Value = 825.0454999999999
$ws.Range("K94").Value = 825.0454999999999
$ws.Range("M94").Value = -374.0454999999999

$ws.Range("H107").Value = 677.5714
$ws.Range("I107").Value = 628.0417
$ws.Range("K107").Value = 628.0417
$ws.Range("M107").Value = 1291.9583

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 169.27272
$ws.Range("I22").Value = 176.2
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 176.2
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 173.8
$ws.Range("N22").Value = -800

$ws.Range("H31").Value = 69802.87
$ws.Range("I31").Value = 113864.89
$ws.Range("K31").Value = 113864.89
$ws.Range("M31").Value = -113569.89

$ws.Range("H34").Value = 69802.87
$ws.Range("I34").Value = 113864.89
$ws.Range("K34").Value = 113864.89
$ws.Range("M34").Value = -113662.89

$ws.Range("H58").Value = 3113.0588
$ws.Range("I58").Value = 2594.6667
$ws.Range("J58").Value = 4357.2
$ws.Range("K58").Value = 2594.6667
$ws.Range("L58").Value = 4357.2
$ws.Range("M58").Value = -2391.6667
$ws.Range("N58").Value = -4763.2

$ws.Range("H136").Value = 3113.0588
$ws.Range("I136").Value = 2594.6667
$ws.Range("J136").Value = 4357.2
$ws.Range("K136").Value = 7784.000100000001
$ws.Range("L136").Value = 13071.6
$ws.Range("M136").Value = -5234.000100000001
$ws.Range("N136").Value = -18171.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 22780.4
$ws.Range("I87").Value = 19686.428
$ws.Range("J87").Value = 29999.666
$ws.Range("K87").Value = 59059.284
$ws.Range("L87").Value = 89998.99800000001
$ws.Range("M87").Value = -57811.284
$ws.Range("N87").Value = -92494.99800000001

$ws.Range("H90").Value = 22780.4
$ws.Range("I90").Value = 19686.428
$ws.Range("J90").Value = 29999.666
$ws.Range("K90").Value = 177177.852
$ws.Range("L90").Value = 269996.994
$ws.Range("M90").Value = -170937.852
$ws.Range("N90").Value = -282476.994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2608.606
$ws.Range("I122").Value = 2237.6155
$ws.Range("J122").Value = 3986.5715
$ws.Range("K122").Value = 6712.8465
$ws.Range("L122").Value = 11959.7145
$ws.Range("M122").Value = -4262.8465
$ws.Range("N122").Value = -16859.7145

$ws.Range("H132").Value = 3036.4546
$ws.Range("I132").Value = 2710.1
$ws.Range("J132").Value = 6300
$ws.Range("K132").Value = 8130.299999999999
$ws.Range("L132").Value = 18900
$ws.Range("M132").Value = -5600.299999999999
$ws.Range("N132").Value = -23960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 127.333336
$ws.Range("I55").Value = 109.5
$ws.Range("K55").Value = 109.5
$ws.Range("M55").Value = 63.5

$ws.Range("H61").Value = 39945.52
$ws.Range("J61").Value = 2421.6
$ws.Range("L61").Value = 2421.6
$ws.Range("N61").Value = -2825.6

$ws.Range("H100").Value = 2448.1538
$ws.Range("I100").Value = 2304.625
$ws.Range("K100").Value = 2304.625
$ws.Range("M100").Value = -1763.625

$ws.Range("H113").Value = 39945.52
$ws.Range("J113").Value = 2421.6
$ws.Range("L113").Value = 2421.6
$ws.Range("N113").Value = -6761.6

$ws.Range("H131").Value = 48498
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 48498
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 48498
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -58578

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1504.9
$ws.Range("I96").Value = 958.1667
$ws.Range("J96").Value = 2325
$ws.Range("K96").Value = 958.1667
$ws.Range("L96").Value = 2325
$ws.Range("M96").Value = 414.8333
$ws.Range("N96").Value = -5071

$ws.Range("H113").Value = 1245.0555
$ws.Range("I113").Value = 1076.2963
$ws.Range("K113").Value = 3228.8889
$ws.Range("M113").Value = -1058.8889

$ws.Range("H132").Value = 2146.2
$ws.Range("I132").Value = 1982.75
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 5948.25
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -3418.25
$ws.Range("N132").Value = -13460
